# Auto-generated edit script: apply the gh-pages data refresh
# (event id=82442 '合肥·全国地下偶像联合公演展-永乐大典Vol.01（取消）' removed,
#  several '想去人数' (want-to-go) counters updated across sheets).
$wb = $excel.ActiveWorkbook

# --- Sheet '展览': refresh '想去人数' counters ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2,6).Value = 3547
$ws.Cells.Item(4,6).Value = 149
$ws.Cells.Item(5,6).Value = 7037
$ws.Cells.Item(6,6).Value = 3600
$ws.Cells.Item(7,6).Value = 63
$ws.Cells.Item(13,6).Value = 36
$ws.Cells.Item(15,6).Value = 601
$ws.Cells.Item(16,6).Value = 54
$ws.Cells.Item(17,6).Value = 49

# --- Sheet '演出': the only listed event was delisted; drop its data row ---
$ws = $wb.Worksheets.Item("演出")
$ws.Rows.Item(2).Delete()

# --- Sheet '本地生活': no changes ---

# --- Sheet '全部类型': refresh counters and remove the delisted event's row,
#     shifting the following events' details up one row ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2,6).Value = 3547
$ws.Cells.Item(3,2).Value = '2024-03-17'
$ws.Cells.Item(3,3).Value = '合肥·CW国潮动漫游戏嘉年华-赵路内场'
$ws.Cells.Item(3,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws.Cells.Item(3,5).Value = '2024.03.17 09:00-03.17 17:00'
$ws.Cells.Item(3,6).Value = 747
$ws.Cells.Item(3,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81954'
$ws.Cells.Item(3,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/2PVn1ahm1708481741272.jpeg'
$ws.Cells.Item(4,2).Value = '2024-03-23'
$ws.Cells.Item(4,3).Value = '合肥·原&铁&崩 only展'
$ws.Cells.Item(4,4).Value = '金寨路与天堂窄路交叉口 梵木艺术中心'
$ws.Cells.Item(4,5).Value = '2024.03.23 09:00-03.23 17:00'
$ws.Cells.Item(4,6).Value = 149
$ws.Cells.Item(4,7).Value = 58
$ws.Cells.Item(4,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81574'
$ws.Cells.Item(4,9).Value = '//i2.hdslb.com/bfs/openplatform/202401/0V5uyX6C1706697212904.png'
$ws.Cells.Item(5,2).Value = '2024-04-04'
$ws.Cells.Item(5,3).Value = '合肥· 第二届漫画城市动漫展 -故事再次开始'
$ws.Cells.Item(5,4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$ws.Cells.Item(5,5).Value = '2024.04.04 09:00-04.05 17:00'
$ws.Cells.Item(5,6).Value = 7037
$ws.Cells.Item(5,7).Value = 60
$ws.Cells.Item(5,8).Value = 'https://show.bilibili.com/platform/detail.html?id=78898'
$ws.Cells.Item(5,9).Value = '//i2.hdslb.com/bfs/openplatform/202402/3NgyB9761708333056023.jpeg'
$ws.Cells.Item(6,3).Value = '合肥·环形宇宙动漫游戏嘉年华'
$ws.Cells.Item(6,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws.Cells.Item(6,5).Value = '2024.04.04 09:30-04.05 17:00'
$ws.Cells.Item(6,6).Value = 3600
$ws.Cells.Item(6,7).Value = 65
$ws.Cells.Item(6,8).Value = 'https://show.bilibili.com/platform/detail.html?id=81916'
$ws.Cells.Item(6,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/1lGzmBT61708336972816.jpeg'
$ws.Cells.Item(7,3).Value = '合肥·环形宇宙动漫游戏嘉年华内场票-谢莹'
$ws.Cells.Item(7,4).Value = '锦绣大道3899号 合肥滨湖会展中心'
$ws.Cells.Item(7,5).Value = '2024.04.04 09:00-04.04 17:00'
$ws.Cells.Item(7,6).Value = 63
$ws.Cells.Item(7,7).Value = 118
$ws.Cells.Item(7,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82515'
$ws.Cells.Item(7,9).Value = '//i2.hdslb.com/bfs/openplatform/202403/L2DFEeao1709800386283.jpeg'
$ws.Cells.Item(8,3).Value = '合肥·环形宇宙动漫游戏嘉年华内场票-钱文青'
$ws.Cells.Item(8,6).Value = 157
$ws.Cells.Item(8,7).Value = 238
$ws.Cells.Item(8,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82375'
$ws.Cells.Item(8,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/3cHtIycW1709692273438.jpeg'
$ws.Cells.Item(9,3).Value = '合肥·第二届漫画城市动漫展内场-柯暮卿'
$ws.Cells.Item(9,4).Value = '凤淮路与固镇路西北角 庐阳全民健身中心'
$ws.Cells.Item(9,5).Value = '2024.04.04 10:00-04.04 17:00'
$ws.Cells.Item(9,6).Value = 30
$ws.Cells.Item(9,7).Value = 158
$ws.Cells.Item(9,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82192'
$ws.Cells.Item(9,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/tcAAj9aj1709193127773.jpeg'
$ws.Cells.Item(10,3).Value = '合肥·第二届漫画城市动漫展内场-风袖'
$ws.Cells.Item(10,6).Value = 39
$ws.Cells.Item(10,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82191'
$ws.Cells.Item(10,9).Value = '//i2.hdslb.com/bfs/openplatform/202402/UZiEzBcc1709192469627.jpeg'
$ws.Cells.Item(11,2).Value = '2024-04-05'
$ws.Cells.Item(11,3).Value = '合肥· 第二届漫画城市动漫展内场-《琅声雅集》'
$ws.Cells.Item(11,5).Value = '2024.04.05 10:00-04.05 17:00'
$ws.Cells.Item(11,6).Value = 90
$ws.Cells.Item(11,7).Value = 168
$ws.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82189'
$ws.Cells.Item(11,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/ZKkx4hTN1709191842946.jpeg'
$ws.Cells.Item(12,2).Value = '2024-04-13'
$ws.Cells.Item(12,3).Value = '合肥·AOO动漫游戏嘉年华'
$ws.Cells.Item(12,4).Value = '芙蓉路291号 正大广场'
$ws.Cells.Item(12,5).Value = '2024.04.13 10:00-04.14 17:00'
$ws.Cells.Item(12,6).Value = 50
$ws.Cells.Item(12,7).Value = 49.9
$ws.Cells.Item(12,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82500'
$ws.Cells.Item(12,9).Value = '//i2.hdslb.com/bfs/openplatform/202403/IpaidRy21709797005042.png'
$ws.Cells.Item(13,2).Value = '2024-04-20'
$ws.Cells.Item(13,3).Value = '合肥·首届运动番only'
$ws.Cells.Item(13,4).Value = '繁华大道6177号 YONEX百胜羽毛球馆(包河店)'
$ws.Cells.Item(13,5).Value = '2024.04.20 10:00-04.20 17:00'
$ws.Cells.Item(13,6).Value = 36
$ws.Cells.Item(13,7).Value = 58
$ws.Cells.Item(13,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82924'
$ws.Cells.Item(13,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/Vw8yXcUF1710489586295.jpeg'
$ws.Cells.Item(14,2).Value = '2024-04-21'
$ws.Cells.Item(14,3).Value = '合肥·银魂only'
$ws.Cells.Item(14,4).Value = '濉溪路118号 合肥栢景假日酒店'
$ws.Cells.Item(14,5).Value = '2024.04.21 09:00-04.21 17:00'
$ws.Cells.Item(14,6).Value = 184
$ws.Cells.Item(14,7).Value = 55
$ws.Cells.Item(14,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82145'
$ws.Cells.Item(14,9).Value = '//i1.hdslb.com/bfs/openplatform/202402/A0Tb5SQ51709091316985.jpeg'
$ws.Cells.Item(15,2).Value = '2024-05-01'
$ws.Cells.Item(15,3).Value = '合肥·Look Look动漫嘉年华'
$ws.Cells.Item(15,4).Value = '新站区东方大道288号 少荃体育中心'
$ws.Cells.Item(15,5).Value = '2024.05.01 10:00-05.01 17:30'
$ws.Cells.Item(15,6).Value = 601
$ws.Cells.Item(15,7).Value = 29.9
$ws.Cells.Item(15,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82311'
$ws.Cells.Item(15,9).Value = '//i2.hdslb.com/bfs/openplatform/202403/jbUNtkAQ1709619599897.png'
$ws.Cells.Item(16,3).Value = '合肥·第十三届次元之门动漫游戏博览会'
$ws.Cells.Item(16,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws.Cells.Item(16,5).Value = '2024.05.01 10:00-05.03 17:00'
$ws.Cells.Item(16,6).Value = 54
$ws.Cells.Item(16,7).Value = '不可售'
$ws.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82910'
$ws.Cells.Item(16,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/NiDA76Th1710471014688.jpeg'
$ws.Cells.Item(17,2).Value = '2024-05-03'
$ws.Cells.Item(17,3).Value = '合肥·BH动漫游戏展'
$ws.Cells.Item(17,4).Value = '科技园路与葡萄园路交口包河区现代农业示范园8号 圩乐田园生态营地'
$ws.Cells.Item(17,5).Value = '2024.05.03 10:00-05.04 16:00'
$ws.Cells.Item(17,6).Value = 49
$ws.Cells.Item(17,7).Value = 40
$ws.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82199'
$ws.Cells.Item(17,9).Value = '//i0.hdslb.com/bfs/openplatform/202402/cSR2xlY61709195356978.jpeg'
$ws.Cells.Item(18,2).Value = '2024-05-18'
$ws.Cells.Item(18,3).Value = '合肥·梦时空SPO1动漫展（取消）'
$ws.Cells.Item(18,4).Value = '阜阳路16号 银瑞林国际大酒店'
$ws.Cells.Item(18,5).Value = '2024.05.18 10:00-05.18 17:00'
$ws.Cells.Item(18,6).Value = 131
$ws.Cells.Item(18,7).Value = '不可售'
$ws.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=80207'
$ws.Cells.Item(18,9).Value = '//i2.hdslb.com/bfs/openplatform/202312/tQQOHYE01703574162111.jpeg'
# the last row is now a duplicate of row 18; remove it to shrink the used range
$ws.Rows.Item(19).Delete()

